$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text representation (no float rounding)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '95.777.70'
$ws.Range("E2").Value = '  -1.57%  '

$ws.Range("D3").Value = '3.627.71'
$ws.Range("E3").Value = '  -2.46%  '

$ws.Range("D4").Value = '2.72'
$ws.Range("E4").Value = '  +40.13%  '

$ws.Range("E5").Value = '  -0.07%  '

$ws.Range("D6").Value = '223.85'
$ws.Range("E6").Value = '  -5.49%  '

$ws.Range("D7").Value = '639.85'
$ws.Range("E7").Value = '  -2.69%  '

$ws.Range("D8").Value = '0.422'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("E9").Value = '  +12.26%  '

$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("D11").Value = '3.625.56'
$ws.Range("E11").Value = '  -2.46%  '

$ws.Range("D12").Value = '48.30'
$ws.Range("E12").Value = '  +7.67%  '

$ws.Range("D13").Value = '0.212'
$ws.Range("E13").Value = '  +1.93%  '

$ws.Range("D14").Value = '0.0000293'
$ws.Range("E14").Value = '  -8.46%  '

$ws.Range("D15").Value = '6.50'
$ws.Range("E15").Value = '  -5.91%  '

$ws.Range("D16").Value = '4.301.45'
$ws.Range("E16").Value = '  -2.54%  '

$ws.Range("D17").Value = '95.517.66'
$ws.Range("E17").Value = '  -1.55%  '

$ws.Range("D18").Value = '22.98'
$ws.Range("E18").Value = '  +22.56%  '

$ws.Range("D19").Value = '8.90'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("D20").Value = '13.91'
$ws.Range("E20").Value = '  +6.62%  '

$ws.Range("D21").Value = '3.625.08'
$ws.Range("E21").Value = '  -2.55%  '

$ws.Range("D22").Value = '0.298'
$ws.Range("E22").Value = '  +54.10%  '

$ws.Range("D23").Value = '0.545'
$ws.Range("E23").Value = '  +7.22%  '

$ws.Range("D24").Value = '516.59'
$ws.Range("E24").Value = '  -1.56%  '

$ws.Range("D25").Value = '3.25'
$ws.Range("E25").Value = '  -6.53%  '

$ws.Range("D26").Value = '125.79'
$ws.Range("E26").Value = '  +18.18%  '

$ws.Range("D27").Value = '0.0000202'
$ws.Range("E27").Value = '  -10.81%  '

$ws.Range("E28").Value = '  -1.32%  '

$ws.Range("D29").Value = '3.811.46'
$ws.Range("E29").Value = '  -2.78%  '

$ws.Range("D30").Value = '12.76'
$ws.Range("E30").Value = '  -5.59%  '

$ws.Range("D31").Value = '13.08'
$ws.Range("E31").Value = '  +3.61%  '

$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("E34").Value = '  +4.86%  '

$ws.Range("E35").Value = '  -5.72%  '

$ws.Range("D36").Value = '32.68'
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("E37").Value = '  +0.28%  '

$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -4.50%  '

$ws.Range("D40").Value = '0.533'
$ws.Range("E40").Value = '  +6.98%  '

$ws.Range("D41").Value = '7.22'
$ws.Range("E41").Value = '  +6.76%  '

$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '584.25'
$ws.Range("E42").Value = '  -9.22%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -4.87%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0521'
$ws.Range("E44").Value = '  +13.90%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '42.08'
$ws.Range("E45").Value = '  +3.34%  '

$ws.Range("E46").Value = '  -0.18%  '

$ws.Range("E47").Value = '  -5.90%  '

$ws.Range("D48").Value = '1.95'
$ws.Range("E48").Value = '  -3.98%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '9.08'
$ws.Range("E49").Value = '  +4.68%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '232.18'
$ws.Range("E50").Value = '  +12.46%  '

$ws.Range("D51").Value = '23.50'
$ws.Range("E51").Value = '  -0.53%  '
